# "Budget incomes for all mun"
#
# The third mini-table on the sheet (rows 17-19, columns B:C) used to be
# the "Investments" block:
#   C17 (header, bold/blue style)   = "Инвестиции"
#   C18 (data row, green style)     = "Инвест. в осн. кап. - invest (тыс. руб) (8109001)"
#   C19 (data row, green style)     = <empty>
#
# It becomes a "Finances" block that also carries a budget-income metric:
#   C17 (header)   = "Финансы"
#   C18 (data row) = "Инвест. в осн. кап. - invest (тыс. руб) (8109001)"  (unchanged)
#   C19 (data row) = "Доходы бюд. - budincome (тыс. руб) (8013001)"      (new)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header of the third table.
$ws.Range("C17").Value = "Финансы"

# C18 keeps the same text (re-asserted for clarity/robustness).
$ws.Range("C18").Value = "Инвест. в осн. кап. - invest (тыс. руб) (8109001)"

# C19 was an empty, formatted cell (same look as C18: fill + border, centered).
# Copy C18's formatting down to C19, then fill in the new budget-income value.
$ws.Range("C18").Copy()
$ws.Range("C19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("C19").Value = "Доходы бюд. - budincome (тыс. руб) (8013001)"

# The author's selection moved from D24 to D26 when the new row was added.
$ws.Range("D26").Select()
